$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-23 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-24 Thursday", 2) | Out-Null
$d.Content.Find.Execute("76÷8=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "37÷2=18, 1", 2) | Out-Null
$d.Content.Find.Execute("26÷9=2, 8", $true, $false, $false, $false, $false, $true, 1, $false, "86÷7=12, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷9=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=22, 0", 2) | Out-Null
$d.Content.Find.Execute("16÷2=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "65÷4=16, 1", 2) | Out-Null
$d.Content.Find.Execute("68÷5=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "83÷5=16, 3", 2) | Out-Null
$d.Content.Find.Execute("74÷5=14, 4", $true, $false, $false, $false, $false, $true, 1, $false, "48÷9=5, 3", 2) | Out-Null
$d.Content.Find.Execute("88÷7=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("96÷4=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "64÷8=8, 0", 2) | Out-Null
$d.Content.Find.Execute("69÷8=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "22÷5=4, 2", 2) | Out-Null
$d.Content.Find.Execute("93÷5=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "54÷2=27, 0", 2) | Out-Null
$d.Content.Find.Execute("61÷3=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=25, 1", 2) | Out-Null
$d.Content.Find.Execute("86÷4=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "29÷4=7, 1", 2) | Out-Null
$d.Content.Find.Execute("85÷7=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=10, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷5=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "38÷7=5, 3", 2) | Out-Null
$d.Content.Find.Execute("45÷3=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=4, 1", 2) | Out-Null
$d.Content.Find.Execute("76÷2=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=10, 0", 2) | Out-Null
$d.Content.Find.Execute("70÷2=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷8=9, 3", 2) | Out-Null
$d.Content.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "15÷2=7, 1", 2) | Out-Null
$d.Content.Find.Execute("51÷5=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=1, 4", 2) | Out-Null
$d.Content.Find.Execute("63÷6=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "66÷9=7, 3", 2) | Out-Null
$d.Content.Find.Execute("34÷6=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "42÷8=5, 2", 2) | Out-Null
$d.Content.Find.Execute("34÷7=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=38, 0", 2) | Out-Null
$d.Content.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 2) | Out-Null
$d.Content.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("63÷8=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "81÷4=20, 1", 2) | Out-Null

Write-Host "Replacements applied."
